# Update the answer key table: each data row (1, 5, 9, 13, 17) of the
# single table holds 5 "two-digit ÷ one-digit" answers. Replace each old
# answer with its new value, cell by cell, using Replace:=1 (wdReplaceOne)
# scoped to that cell's Range so that duplicate answer text elsewhere in
# the table is not accidentally overwritten as well.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=7, 0", 1) | Out-Null

$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("77÷7=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "99÷7=14, 1", 1) | Out-Null

$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("51÷8=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "35÷3=11, 2", 1) | Out-Null

$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("23÷5=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "16÷3=5, 1", 1) | Out-Null

$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("46÷8=5, 6", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=19, 3", 1) | Out-Null

$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("27÷6=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=5, 0", 1) | Out-Null

$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("77÷8=9, 5", $true, $false, $false, $false, $false, $true, 1, $false, "65÷5=13, 0", 1) | Out-Null

$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("35÷4=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "68÷9=7, 5", 1) | Out-Null

$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("26÷7=3, 5", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 1) | Out-Null

$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("76÷7=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "97÷3=32, 1", 1) | Out-Null

$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("74÷9=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "12÷7=1, 5", 1) | Out-Null

$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("77÷2=38, 1", $true, $false, $false, $false, $false, $true, 1, $false, "19÷2=9, 1", 1) | Out-Null

$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("16÷6=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "48÷4=12, 0", 1) | Out-Null

$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("35÷2=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷7=6, 4", 1) | Out-Null

$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("76÷9=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "79÷5=15, 4", 1) | Out-Null

$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("84÷4=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "42÷5=8, 2", 1) | Out-Null

$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("60÷9=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "45÷9=5, 0", 1) | Out-Null

$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("32÷3=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "95÷7=13, 4", 1) | Out-Null

$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("47÷4=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "29÷7=4, 1", 1) | Out-Null

$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("86÷8=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "28÷7=4, 0", 1) | Out-Null

$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("51÷8=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "11÷7=1, 4", 1) | Out-Null

$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("97÷8=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=6, 4", 1) | Out-Null

$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("81÷7=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=15, 1", 1) | Out-Null

$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("54÷8=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "53÷6=8, 5", 1) | Out-Null

$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("57÷5=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=5, 5", 1) | Out-Null
